# "Taking latest changes and appending my changes"
# Update the "Test Cases" sheet (sheet1): rows 7 & 8 Results flip from PASS -> SKIP,
# and four brand-new rows (9-12) are appended describing new OPQA-215..218 test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Rows 7 & 8: Results column flips PASS -> SKIP, and the Runmode cell's
#     format is normalized to match the rest of the Runmode column ---
$ws.Range("E7").Value = "SKIP"
$ws.Range("E8").Value = "SKIP"

$ws.Range("D3").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = "Y"

$ws.Range("D3").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = "Y"

# --- New rows 9-12 ---
# Seed formatting for each new row by pasting the formats of an existing,
# similarly-styled row, then overwrite the values on top.

# Row 9  (style-wise a sibling of rows 5/6: A/B bold-ish "7" style, C/D plain, E plain)
$ws.Range("A5:E5").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A9").Value = "TestCase_F8"
$ws.Range("B9").Value = "OPQA-215"
$ws.Range("C9").Value = "Verify that user able to recevies a notification when other user commented on his post"
$ws.Range("D9").Value = "Y"
$ws.Range("E9").Value = "SKIP"

# Row 10 (same family as row 9)
$ws.Range("A5:E5").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A10").Value = "TestCase_F9"
$ws.Range("B10").Value = "OPQA-216"
$ws.Range("C10").Value = "Verify that user receives a notification when someone he is following user comments on a post"
$ws.Range("D10").Value = "Y"
$ws.Range("E10").Value = "SKIP"

# Row 11 (same family as row 9, but Results = PASS)
$ws.Range("A5:E5").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A11").Value = "TestCase_F10"
$ws.Range("B11").Value = "OPQA-217"
$ws.Range("C11").Value = "Verify that user receives a notification when someone comments on an post contained in his watchlist"
$ws.Range("D11").Value = "Y"
$ws.Range("E11").Value = "PASS"

# Row 12 (style-wise a sibling of rows 2-4: plain "3" style on A, PASS result)
$ws.Range("A4:E4").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A12").Value = "TestCase_F11"
$ws.Range("B12").Value = "OPQA-218"
$ws.Range("C12").Value = "Verify that user receives a notification if someone likes his comment on a post"
$ws.Range("D12").Value = "Y"
$ws.Range("E12").Value = "PASS"

# --- View bookkeeping: selection moves to D9, used range now spans through row 12 ---
$ws.Range("D9").Select() | Out-Null
